$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1252.35
$ws.Range("I8").Value = 156.75
$ws.Range("J8").Value = 1526.25
$ws.Range("K8").Value = 470.25
$ws.Range("L8").Value = 4578.75
$ws.Range("M8").Value = -331.25
$ws.Range("N8").Value = -4856.75

$ws.Range("H9").Value = 576.9231
$ws.Range("I9").Value = 437.5
$ws.Range("J9").Value = 800
$ws.Range("K9").Value = 437.5
$ws.Range("L9").Value = 800
$ws.Range("M9").Value = -268.5
$ws.Range("N9").Value = -1138

$ws.Range("H98").Value = 2025.4736
$ws.Range("I98").Value = 1895.7778
$ws.Range("J98").Value = 2142.2
$ws.Range("K98").Value = 1895.7778
$ws.Range("L98").Value = 2142.2
$ws.Range("M98").Value = -397.7778000000001
$ws.Range("N98").Value = -5138.2

$ws.Range("H112").Value = 3969414.2
$ws.Range("I112").Value = 1850
$ws.Range("J112").Value = 4546514.5
$ws.Range("K112").Value = 5550
$ws.Range("L112").Value = 13639543.5
$ws.Range("M112").Value = -4442
$ws.Range("N112").Value = -13641759.5

$ws.Range("H122").Value = 2025.4736
$ws.Range("I122").Value = 1895.7778
$ws.Range("J122").Value = 2142.2
$ws.Range("K122").Value = 5687.3334
$ws.Range("L122").Value = 6426.599999999999
$ws.Range("M122").Value = -3237.3334
$ws.Range("N122").Value = -11326.6

$ws.Range("H133").Value = 46250
$ws.Range("J133").Value = 46250
$ws.Range("L133").Value = 46250
$ws.Range("N133").Value = -56370

$ws.Range("H137").Value = 3114.2285
$ws.Range("I137").Value = 2790.2903
$ws.Range("J137").Value = 5624.75
$ws.Range("K137").Value = 8370.8709
$ws.Range("L137").Value = 16874.25
$ws.Range("M137").Value = -5820.8709
$ws.Range("N137").Value = -21974.25

$ws.Range("H138").Value = 3114.2588
$ws.Range("I138").Value = 1775.9722
$ws.Range("J138").Value = 4097.4897
$ws.Range("K138").Value = 5327.9166
$ws.Range("L138").Value = 12292.4691
$ws.Range("M138").Value = -187.9165999999996
$ws.Range("N138").Value = -22572.4691

$ws.Range("H141").Value = 317922.6
$ws.Range("I141").Value = 1037.25
$ws.Range("K141").Value = 3111.75
$ws.Range("M141").Value = 2068.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2503.4167
$ws.Range("I61").Value = 1276.8636
$ws.Range("J61").Value = 4430.857
$ws.Range("K61").Value = 1276.8636
$ws.Range("L61").Value = 4430.857
$ws.Range("M61").Value = -1064.8636
$ws.Range("N61").Value = -4854.857

$ws.Range("H132").Value = 22730424
$ws.Range("I132").Value = 50002250
$ws.Range("J132").Value = 3905.375
$ws.Range("K132").Value = 150006750
$ws.Range("L132").Value = 11716.125
$ws.Range("M132").Value = -150004220
$ws.Range("N132").Value = -16776.125

$ws.Range("H133").Value = 29660
$ws.Range("J133").Value = 29660
$ws.Range("L133").Value = 29660
$ws.Range("N133").Value = -34720

$ws.Range("H136").Value = 2503.4167
$ws.Range("I136").Value = 1276.8636
$ws.Range("J136").Value = 4430.857
$ws.Range("K136").Value = 3830.5908
$ws.Range("L136").Value = 13292.571
$ws.Range("M136").Value = -1280.5908
$ws.Range("N136").Value = -18392.571

$ws.Range("H141").Value = 29471.428
$ws.Range("J141").Value = 29471.428
$ws.Range("L141").Value = 29471.428
$ws.Range("N141").Value = -39831.428

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").Value = ""

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").Value = ""

$ws.Range("H75").Value = 22553.5
$ws.Range("I75").Value = 5107
$ws.Range("K75").Value = 5107
$ws.Range("M75").Value = -4171

$ws.Range("H78").Value = 22553.5
$ws.Range("I78").Value = 5107
$ws.Range("K78").Value = 15321
$ws.Range("M78").Value = -10641

$ws.Range("H134").Value = 2702.3333
$ws.Range("I134").Value = 1780.1538
$ws.Range("J134").Value = 5100
$ws.Range("K134").Value = 5340.4614
$ws.Range("L134").Value = 15300
$ws.Range("M134").Value = -2805.4614
$ws.Range("N134").Value = -20370

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 13000
$ws.Range("I56").Value = 9000
$ws.Range("J56").Value = 15000
$ws.Range("K56").Value = 9000
$ws.Range("L56").Value = 15000
$ws.Range("M56").Value = -8155
$ws.Range("N56").Value = -16690

$ws.Range("H132").Value = 2475.25
$ws.Range("I132").Value = 1951.5555
$ws.Range("J132").Value = 3562.923
$ws.Range("K132").Value = 5854.666499999999
$ws.Range("L132").Value = 10688.769
$ws.Range("M132").Value = -3324.666499999999
$ws.Range("N132").Value = -15748.769

$ws.Range("H141").Value = 30000
$ws.Range("J141").Value = 30000
$ws.Range("L141").Value = 30000
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1122.55
$ws.Range("I5").Value = 636.1667
$ws.Range("J5").Value = 5500
$ws.Range("K5").Value = 1908.5001
$ws.Range("L5").Value = 16500
$ws.Range("M5").Value = -1796.5001
$ws.Range("N5").Value = -16724

$ws.Range("H122").Value = 1773.25
$ws.Range("I122").Value = 596.25
$ws.Range("J122").Value = 2361.75
$ws.Range("K122").Value = 5366.25
$ws.Range("L122").Value = 21255.75
$ws.Range("M122").Value = -2916.25
$ws.Range("N122").Value = -26155.75

$ws.Range("H132").Value = 1929
$ws.Range("J132").Value = 2708.5715
$ws.Range("L132").Value = 24377.1435
$ws.Range("N132").Value = -29437.1435

$ws.Range("H135").Value = 1122.55
$ws.Range("I135").Value = 636.1667
$ws.Range("J135").Value = 5500
$ws.Range("K135").Value = 5725.5003
$ws.Range("L135").Value = 49500
$ws.Range("M135").Value = -3190.5003
$ws.Range("N135").Value = -54570

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6552.6113
$ws.Range("I122").Value = 5924.8335
$ws.Range("K122").Value = 17774.5005
$ws.Range("M122").Value = -15324.5005

$ws.Range("H132").Value = 3440.4243
$ws.Range("I132").Value = 2965.0527
$ws.Range("J132").Value = 4085.5715
$ws.Range("K132").Value = 8895.158100000001
$ws.Range("L132").Value = 12256.7145
$ws.Range("M132").Value = -6365.158100000001
$ws.Range("N132").Value = -17316.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H57").Value = 2608
$ws.Range("I57").Value = 2010
$ws.Range("K57").Value = 2010
$ws.Range("M57").Value = -1444

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""

$ws.Range("H136").Value = 2545.3428
$ws.Range("I136").Value = 2210.6296
$ws.Range("J136").Value = 3675
$ws.Range("K136").Value = 6631.888800000001
$ws.Range("L136").Value = 11025
$ws.Range("M136").Value = -4081.888800000001
$ws.Range("N136").Value = -16125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4227.1836
$ws.Range("I132").Value = 1780.3928
$ws.Range("J132").Value = 7489.5713
$ws.Range("K132").Value = 5341.178400000001
$ws.Range("L132").Value = 22468.7139
$ws.Range("M132").Value = -2811.178400000001
$ws.Range("N132").Value = -27528.7139
